$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title heading.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $metaPara = $d.Paragraphs.Item($i)
    if ($metaPara.Range.Text.StartsWith("Meta description")) {
        $metaPara.Range.Delete()
        break
    }
}

# 2. Insert a new bold "Play Auspicious Fortune God for Free - Review 2021" paragraph
#    right before the final (image-prompt) paragraph.
$imgParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Create an eye-catching feature image")) {
        $imgParaIndex = $i
        break
    }
}
$imgPara = $d.Paragraphs.Item($imgParaIndex)
$imgPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($imgParaIndex)
$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Auspicious Fortune God for Free - Review 2021</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newXml)

# 3. Replace the final paragraph's (italic) image-prompt text with the meta description text.
$finalPara = $d.Paragraphs.Item($imgParaIndex + 1)
$oldText = 'Create an eye-catching feature image for "Auspicious Fortune God" that showcases a happy Maya warrior with glasses in cartoon style. Use bright and bold colors that are reminiscent of the game''s Chinese-theme, such as red and gold, to catch the viewer''s attention. Position the warrior in a confident and charismatic pose, with his hands on his hips and a big smile on his face. Make sure to incorporate the game''s logo into the image and any other relevant symbols such as fortune deities and money trees. The overall image should convey a sense of excitement and fun while highlighting the game''s unique features.'
$newText = "Explore the features of Auspicious Fortune God. Read our unbiased review and play for free. Discover volatility, betting range and responsible gambling advice."
$finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
